$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 91, shifting existing rows 91-99 down to 92-100
$ws.Rows.Item(91).Insert()

# Populate the newly inserted row 91 with the new weekly record
$ws.Range("A91").Value = 4
$ws.Range("B91").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C91").Value = "Los Lagos"
$ws.Range("D91").Value = 45258
$ws.Range("E91").Value = 10
$ws.Range("F91").Value = 100112012
$ws.Range("G91").Value = "Espinaca"
$ws.Range("H91").Value = "Sin especificar"
$ws.Range("I91").Value = "Primera"
$ws.Range("J91").Value = 35
$ws.Range("K91").Value = 20000
$ws.Range("L91").Value = 20000
$ws.Range("M91").Value = 20000
$ws.Range("N91").Value = '$/cuna 10 kilos'
$ws.Range("O91").Value = "Región Metropolitana"
$ws.Range("P91").Value = 2000
$ws.Range("Q91").Value = 10
$ws.Range("R91").Value = "Hortaliza"
